$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "41.891.84"
$ws.Range('E2').Value = "  +4.86%  "
$ws.Range('D3').Value = "2.267.09"
$ws.Range('E3').Value = "  +2.19%  "
$ws.Range('E4').Value = "  -0.03%  "
$ws.Range('D5').Value = "'301.85"
$ws.Range('E5').Value = "  +3.40%  "
$ws.Range('D6').Value = "'92.05"
$ws.Range('E6').Value = "  +6.09%  "
$ws.Range('D7').Value = "'0.531"
$ws.Range('E7').Value = "  +3.35%  "
$ws.Range('E8').Value = "  -0.01%  "
$ws.Range('E9').Value = "  +4.28%  "
$ws.Range('D10').Value = "'54.39"
$ws.Range('E10').Value = "  +7.98%  "
$ws.Range('D11').Value = "'32.26"
$ws.Range('E11').Value = "  +6.00%  "
$ws.Range('D12').Value = "'0.0798"
$ws.Range('E12').Value = "  +2.45%  "
$ws.Range('E14').Value = "  +3.62%  "
$ws.Range('D15').Value = "2.616.73"
$ws.Range('E15').Value = "  +2.07%  "
$ws.Range('D16').Value = "'14.14"
$ws.Range('E16').Value = "  +2.76%  "
$ws.Range('D17').Value = "2.268.69"
$ws.Range('E17').Value = "  +0.94%  "
$ws.Range('E18').Value = "  +3.64%  "
$ws.Range('D19').Value = "41.793.31"
$ws.Range('E19').Value = "  +4.80%  "
$ws.Range('D20').Value = "'12.12"
$ws.Range('E20').Value = "  +8.96%  "
$ws.Range('D21').Value = "0.0₃0904"
$ws.Range('E21').Value = "  +2.03%  "
$ws.Range('D22').Value = "'5.94"
$ws.Range('E22').Value = "  +3.36%  "
$ws.Range('D23').Value = "'66.99"
$ws.Range('E23').Value = "  +2.24%  "
$ws.Range('D24').Value = "'241.71"
$ws.Range('E24').Value = "  +1.98%  "
$ws.Range('D25').Value = "'2.57"
$ws.Range('E25').Value = "  +4.10%  "
$ws.Range('E26').Value = "  -0.03%  "
$ws.Range('E27').Value = "  +3.98%  "
$ws.Range('D28').Value = "'23.92"
$ws.Range('E28').Value = "  +2.79%  "
$ws.Range('E29').Value = "  +4.57%  "
$ws.Range('E30').Value = "  -5.45%  "
$ws.Range('D31').Value = "'159.33"
$ws.Range('E31').Value = "  +1.19%  "
$ws.Range('D32').Value = "'33.79"
$ws.Range('E32').Value = "  +6.33%  "
$ws.Range('E33').Value = "  -0.05%  "
$ws.Range('E34').Value = "  +4.20%  "
$ws.Range('E35').Value = "  +4.43%  "
$ws.Range('D36').Value = "'3.07"
$ws.Range('E36').Value = "  +2.89%  "
$ws.Range('E37').Value = "  +2.10%  "
$ws.Range('B38').Value = "Kaspa"
$ws.Range('C38').Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D38').Value = "'0.104"
$ws.Range('E38').Value = "  +5.11%  "
$ws.Range('B39').Value = "Celestia"
$ws.Range('C39').Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range('D39').Value = "'16.55"
$ws.Range('E39').Value = "  +8.98%  "
$ws.Range('D40').Value = "'0.115"
$ws.Range('E40').Value = "  +3.31%  "
$ws.Range('E41').Value = "  +5.21%  "
$ws.Range('E42').Value = "  +6.00%  "
$ws.Range('D43').Value = "2.069.63"
$ws.Range('E43').Value = "  -0.82%  "
$ws.Range('D44').Value = "'19.86"
$ws.Range('E44').Value = "  +10.52%  "
$ws.Range('E45').Value = "  +3.45%  "
$ws.Range('D46').Value = "'10.21"
$ws.Range('E46').Value = "  +4.42%  "
$ws.Range('D47').Value = "'2.91"
$ws.Range('E47').Value = "  +7.84%  "
$ws.Range('E48').Value = "  +2.48%  "
$ws.Range('E49').Value = "  +3.84%  "
$ws.Range('D50').Value = "'1.15"
$ws.Range('E50').Value = "  +3.57%  "
$ws.Range('D51').Value = "'51.74"
$ws.Range('E51').Value = "  +5.62%  "
